$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The naive-forecaster component table shifts one column to the left (oldest
# forecast-origin column dropped) and the two now-empty trailing rows go away.
# Clear every cell first (old layout had data out to column BA / row 24) so no
# stale values survive from cells that fall outside the new footprint, then
# re-populate the surviving cells with their corrected values.
$ws.Range("A1:BA24").ClearContents() | Out-Null
$ws.Range("A23:A24").EntireRow.Delete() | Out-Null
$ws.Range("BA1:BA1").EntireColumn.Delete() | Out-Null

# Column A: period end dates (unchanged)
$colA = New-Object 'object[,]' 21,1
$colA[0,0] = 39813
$colA[1,0] = 40178
$colA[2,0] = 40543
$colA[3,0] = 40908
$colA[4,0] = 41274
$colA[5,0] = 41639
$colA[6,0] = 42004
$colA[7,0] = 42369
$colA[8,0] = 42735
$colA[9,0] = 43100
$colA[10,0] = 43465
$colA[11,0] = 43830
$colA[12,0] = 44196
$colA[13,0] = 44561
$colA[14,0] = 44926
$colA[15,0] = 45291
$colA[16,0] = 45657
$colA[17,0] = 46022
$colA[18,0] = 46387
$colA[19,0] = 46752
$colA[20,0] = 47118
$ws.Range("A2:A22").Value = $colA

# Row 1 (forecast-origin header dates) and rows 3-21 (forecast component values)
$row1 = New-Object 'object[,]' 1,51
$row1[0,0] = 39583
$row1[0,1] = 39765
$row1[0,2] = 39948
$row1[0,3] = 40130
$row1[0,4] = 40310
$row1[0,5] = 40494
$row1[0,6] = 40676
$row1[0,7] = 40862
$row1[0,8] = 41044
$row1[0,9] = 41228
$row1[0,10] = 41409
$row1[0,11] = 41592
$row1[0,12] = 41774
$row1[0,13] = 41957
$row1[0,14] = 42137
$row1[0,15] = 42321
$row1[0,16] = 42503
$row1[0,17] = 42689
$row1[0,18] = 42867
$row1[0,19] = 43053
$row1[0,20] = 43145
$row1[0,21] = 43235
$row1[0,22] = 43326
$row1[0,23] = 43418
$row1[0,24] = 43510
$row1[0,25] = 43600
$row1[0,26] = 43691
$row1[0,27] = 43783
$row1[0,28] = 43875
$row1[0,29] = 43966
$row1[0,30] = 44068
$row1[0,31] = 44159
$row1[0,32] = 44251
$row1[0,33] = 44341
$row1[0,34] = 44432
$row1[0,35] = 44525
$row1[0,36] = 44617
$row1[0,37] = 44706
$row1[0,38] = 44798
$row1[0,39] = 44890
$row1[0,40] = 44981
$row1[0,41] = 45071
$row1[0,42] = 45163
$row1[0,43] = 45254
$row1[0,44] = 45345
$row1[0,45] = 45436
$row1[0,46] = 45534
$row1[0,47] = 45618
$row1[0,48] = 45713
$row1[0,49] = 45800
$row1[0,50] = 45891
$ws.Range("B1:AZ1").Value = $row1

$row3 = New-Object 'object[,]' 1,51
$row3[0,0] = 6.992417114397731
$row3[0,1] = 2.651903832953884
$row3[0,2] = 1.052572528399653
$row3[0,3] = 0.1715429114845346
$row3[0,4] = 0.1715429114845346
$row3[0,5] = 0.1715429114845346
$row3[0,6] = 0.1715429114845346
$row3[0,7] = 0.1715429114845346
$row3[0,8] = 0.1715429114845346
$row3[0,9] = 0.1715429114845346
$row3[0,10] = 0.1715429114845346
$row3[0,11] = 0.1715429114845346
$row3[0,12] = 0.1715429114845346
$row3[0,13] = 0.1715429114845346
$row3[0,14] = 0.1715429114845346
$row3[0,15] = 0.1715429114845346
$row3[0,16] = 0.1715429114845346
$row3[0,17] = 0.1715429114845346
$row3[0,18] = 0.1715429114845346
$row3[0,19] = 0.1715429114845346
$row3[0,20] = 0.1715429114845346
$row3[0,21] = 0.1715429114845346
$row3[0,22] = 0.1715429114845346
$row3[0,23] = 0.1715429114845346
$row3[0,24] = 0.1715429114845346
$row3[0,25] = 0.1715429114845346
$row3[0,26] = 0.1715429114845346
$row3[0,27] = 0.1715429114845346
$row3[0,28] = 0.1715429114845346
$row3[0,29] = 0.1715429114845346
$row3[0,30] = 0.1715429114845346
$row3[0,31] = 0.1715429114845346
$row3[0,32] = 0.1715429114845346
$row3[0,33] = 0.1715429114845346
$row3[0,34] = 0.1715429114845346
$row3[0,35] = 0.1715429114845346
$row3[0,36] = 0.1715429114845346
$row3[0,37] = 0.1715429114845346
$row3[0,38] = 0.1715429114845346
$row3[0,39] = 0.1715429114845346
$row3[0,40] = 0.1715429114845346
$row3[0,41] = 0.1715429114845346
$row3[0,42] = 0.1715429114845346
$row3[0,43] = 0.1715429114845346
$row3[0,44] = 0.1715429114845346
$row3[0,45] = 0.1715429114845346
$row3[0,46] = 0.1715429114845346
$row3[0,47] = 0.1715429114845346
$row3[0,48] = 0.1715429114845346
$row3[0,49] = 0.1715429114845346
$row3[0,50] = 0.1715429114845346
$ws.Range("B3:AZ3").Value = $row3

$row4 = New-Object 'object[,]' 1,51
$row4[0,0] = 7.18590312890619
$row4[0,1] = 3.648892256099967
$row4[0,2] = 2.396905458966625
$row4[0,3] = 0.5356365903280924
$row4[0,4] = 0.475547144596522
$row4[0,5] = 0.8004663283405655
$row4[0,6] = 0.8004663283405655
$row4[0,7] = 0.8004663283405655
$row4[0,8] = 0.8004663283405655
$row4[0,9] = 0.8004663283405655
$row4[0,10] = 0.8004663283405655
$row4[0,11] = 0.8004663283405655
$row4[0,12] = 0.8004663283405655
$row4[0,13] = 0.8004663283405655
$row4[0,14] = 0.8004663283405655
$row4[0,15] = 0.8004663283405655
$row4[0,16] = 0.8004663283405655
$row4[0,17] = 0.8004663283405655
$row4[0,18] = 0.8004663283405655
$row4[0,19] = 0.8004663283405655
$row4[0,20] = 0.8004663283405655
$row4[0,21] = 0.8004663283405655
$row4[0,22] = 0.8004663283405655
$row4[0,23] = 0.8004663283405655
$row4[0,24] = 0.8004663283405655
$row4[0,25] = 0.8004663283405655
$row4[0,26] = 0.8004663283405655
$row4[0,27] = 0.8004663283405655
$row4[0,28] = 0.8004663283405655
$row4[0,29] = 0.8004663283405655
$row4[0,30] = 0.8004663283405655
$row4[0,31] = 0.8004663283405655
$row4[0,32] = 0.8004663283405655
$row4[0,33] = 0.8004663283405655
$row4[0,34] = 0.8004663283405655
$row4[0,35] = 0.8004663283405655
$row4[0,36] = 0.8004663283405655
$row4[0,37] = 0.8004663283405655
$row4[0,38] = 0.8004663283405655
$row4[0,39] = 0.8004663283405655
$row4[0,40] = 0.8004663283405655
$row4[0,41] = 0.8004663283405655
$row4[0,42] = 0.8004663283405655
$row4[0,43] = 0.8004663283405655
$row4[0,44] = 0.8004663283405655
$row4[0,45] = 0.8004663283405655
$row4[0,46] = 0.8004663283405655
$row4[0,47] = 0.8004663283405655
$row4[0,48] = 0.8004663283405655
$row4[0,49] = 0.8004663283405655
$row4[0,50] = 0.8004663283405655
$ws.Range("B4:AZ4").Value = $row4

$row5 = New-Object 'object[,]' 1,49
$row5[0,0] = 2.557501775704929
$row5[0,1] = 1.508458613525354
$row5[0,2] = 1.495774185788745
$row5[0,3] = 1.805987045940682
$row5[0,4] = 4.312104569761943
$row5[0,5] = 5.253783907501819
$row5[0,6] = 5.253783907501819
$row5[0,7] = 5.253783907501819
$row5[0,8] = 5.253783907501819
$row5[0,9] = 5.253783907501819
$row5[0,10] = 5.253783907501819
$row5[0,11] = 5.253783907501819
$row5[0,12] = 5.253783907501819
$row5[0,13] = 5.253783907501819
$row5[0,14] = 5.253783907501819
$row5[0,15] = 5.253783907501819
$row5[0,16] = 5.253783907501819
$row5[0,17] = 5.253783907501819
$row5[0,18] = 5.253783907501819
$row5[0,19] = 5.253783907501819
$row5[0,20] = 5.253783907501819
$row5[0,21] = 5.253783907501819
$row5[0,22] = 5.253783907501819
$row5[0,23] = 5.253783907501819
$row5[0,24] = 5.253783907501819
$row5[0,25] = 5.253783907501819
$row5[0,26] = 5.253783907501819
$row5[0,27] = 5.253783907501819
$row5[0,28] = 5.253783907501819
$row5[0,29] = 5.253783907501819
$row5[0,30] = 5.253783907501819
$row5[0,31] = 5.253783907501819
$row5[0,32] = 5.253783907501819
$row5[0,33] = 5.253783907501819
$row5[0,34] = 5.253783907501819
$row5[0,35] = 5.253783907501819
$row5[0,36] = 5.253783907501819
$row5[0,37] = 5.253783907501819
$row5[0,38] = 5.253783907501819
$row5[0,39] = 5.253783907501819
$row5[0,40] = 5.253783907501819
$row5[0,41] = 5.253783907501819
$row5[0,42] = 5.253783907501819
$row5[0,43] = 5.253783907501819
$row5[0,44] = 5.253783907501819
$row5[0,45] = 5.253783907501819
$row5[0,46] = 5.253783907501819
$row5[0,47] = 5.253783907501819
$row5[0,48] = 5.253783907501819
$ws.Range("D5:AZ5").Value = $row5

$row6 = New-Object 'object[,]' 1,47
$row6[0,0] = 1.488234279941647
$row6[0,1] = 0.8024032015999882
$row6[0,2] = 1.976172106438545
$row6[0,3] = 4.382531137514767
$row6[0,4] = 4.639864760432189
$row6[0,5] = 3.522405026196918
$row6[0,6] = 3.522405026196918
$row6[0,7] = 3.522405026196918
$row6[0,8] = 3.522405026196918
$row6[0,9] = 3.522405026196918
$row6[0,10] = 3.522405026196918
$row6[0,11] = 3.522405026196918
$row6[0,12] = 3.522405026196918
$row6[0,13] = 3.522405026196918
$row6[0,14] = 3.522405026196918
$row6[0,15] = 3.522405026196918
$row6[0,16] = 3.522405026196918
$row6[0,17] = 3.522405026196918
$row6[0,18] = 3.522405026196918
$row6[0,19] = 3.522405026196918
$row6[0,20] = 3.522405026196918
$row6[0,21] = 3.522405026196918
$row6[0,22] = 3.522405026196918
$row6[0,23] = 3.522405026196918
$row6[0,24] = 3.522405026196918
$row6[0,25] = 3.522405026196918
$row6[0,26] = 3.522405026196918
$row6[0,27] = 3.522405026196918
$row6[0,28] = 3.522405026196918
$row6[0,29] = 3.522405026196918
$row6[0,30] = 3.522405026196918
$row6[0,31] = 3.522405026196918
$row6[0,32] = 3.522405026196918
$row6[0,33] = 3.522405026196918
$row6[0,34] = 3.522405026196918
$row6[0,35] = 3.522405026196918
$row6[0,36] = 3.522405026196918
$row6[0,37] = 3.522405026196918
$row6[0,38] = 3.522405026196918
$row6[0,39] = 3.522405026196918
$row6[0,40] = 3.522405026196918
$row6[0,41] = 3.522405026196918
$row6[0,42] = 3.522405026196918
$row6[0,43] = 3.522405026196918
$row6[0,44] = 3.522405026196918
$row6[0,45] = 3.522405026196918
$row6[0,46] = 3.522405026196918
$ws.Range("F6:AZ6").Value = $row6

$row7 = New-Object 'object[,]' 1,45
$row7[0,0] = 1.893295577996779
$row7[0,1] = 3.197643605100797
$row7[0,2] = 4.074582884048139
$row7[0,3] = 1.988448192515935
$row7[0,4] = 1.418316910291906
$row7[0,5] = 1.656063945467268
$row7[0,6] = 1.656063945467268
$row7[0,7] = 1.656063945467268
$row7[0,8] = 1.656063945467268
$row7[0,9] = 1.656063945467268
$row7[0,10] = 1.656063945467268
$row7[0,11] = 1.656063945467268
$row7[0,12] = 1.656063945467268
$row7[0,13] = 1.656063945467268
$row7[0,14] = 1.656063945467268
$row7[0,15] = 1.656063945467268
$row7[0,16] = 1.656063945467268
$row7[0,17] = 1.656063945467268
$row7[0,18] = 1.656063945467268
$row7[0,19] = 1.656063945467268
$row7[0,20] = 1.656063945467268
$row7[0,21] = 1.656063945467268
$row7[0,22] = 1.656063945467268
$row7[0,23] = 1.656063945467268
$row7[0,24] = 1.656063945467268
$row7[0,25] = 1.656063945467268
$row7[0,26] = 1.656063945467268
$row7[0,27] = 1.656063945467268
$row7[0,28] = 1.656063945467268
$row7[0,29] = 1.656063945467268
$row7[0,30] = 1.656063945467268
$row7[0,31] = 1.656063945467268
$row7[0,32] = 1.656063945467268
$row7[0,33] = 1.656063945467268
$row7[0,34] = 1.656063945467268
$row7[0,35] = 1.656063945467268
$row7[0,36] = 1.656063945467268
$row7[0,37] = 1.656063945467268
$row7[0,38] = 1.656063945467268
$row7[0,39] = 1.656063945467268
$row7[0,40] = 1.656063945467268
$row7[0,41] = 1.656063945467268
$row7[0,42] = 1.656063945467268
$row7[0,43] = 1.656063945467268
$row7[0,44] = 1.656063945467268
$ws.Range("H7:AZ7").Value = $row7

$row8 = New-Object 'object[,]' 1,43
$row8[0,0] = 4.184092216308799
$row8[0,1] = 3.608060659590451
$row8[0,2] = 3.247860853607465
$row8[0,3] = 3.166945525867848
$row8[0,4] = 3.520945360626571
$row8[0,5] = 4.06235252733802
$row8[0,6] = 4.06235252733802
$row8[0,7] = 4.06235252733802
$row8[0,8] = 4.06235252733802
$row8[0,9] = 4.06235252733802
$row8[0,10] = 4.06235252733802
$row8[0,11] = 4.06235252733802
$row8[0,12] = 4.06235252733802
$row8[0,13] = 4.06235252733802
$row8[0,14] = 4.06235252733802
$row8[0,15] = 4.06235252733802
$row8[0,16] = 4.06235252733802
$row8[0,17] = 4.06235252733802
$row8[0,18] = 4.06235252733802
$row8[0,19] = 4.06235252733802
$row8[0,20] = 4.06235252733802
$row8[0,21] = 4.06235252733802
$row8[0,22] = 4.06235252733802
$row8[0,23] = 4.06235252733802
$row8[0,24] = 4.06235252733802
$row8[0,25] = 4.06235252733802
$row8[0,26] = 4.06235252733802
$row8[0,27] = 4.06235252733802
$row8[0,28] = 4.06235252733802
$row8[0,29] = 4.06235252733802
$row8[0,30] = 4.06235252733802
$row8[0,31] = 4.06235252733802
$row8[0,32] = 4.06235252733802
$row8[0,33] = 4.06235252733802
$row8[0,34] = 4.06235252733802
$row8[0,35] = 4.06235252733802
$row8[0,36] = 4.06235252733802
$row8[0,37] = 4.06235252733802
$row8[0,38] = 4.06235252733802
$row8[0,39] = 4.06235252733802
$row8[0,40] = 4.06235252733802
$row8[0,41] = 4.06235252733802
$row8[0,42] = 4.06235252733802
$ws.Range("J8:AZ8").Value = $row8

$row9 = New-Object 'object[,]' 1,41
$row9[0,0] = 3.343464100342031
$row9[0,1] = 3.019459040387984
$row9[0,2] = 2.547371915279606
$row9[0,3] = 3.684750195712683
$row9[0,4] = 3.285568146716344
$row9[0,5] = 3.05427116350534
$row9[0,6] = 3.05427116350534
$row9[0,7] = 3.05427116350534
$row9[0,8] = 3.05427116350534
$row9[0,9] = 3.05427116350534
$row9[0,10] = 3.05427116350534
$row9[0,11] = 3.05427116350534
$row9[0,12] = 3.05427116350534
$row9[0,13] = 3.05427116350534
$row9[0,14] = 3.05427116350534
$row9[0,15] = 3.05427116350534
$row9[0,16] = 3.05427116350534
$row9[0,17] = 3.05427116350534
$row9[0,18] = 3.05427116350534
$row9[0,19] = 3.05427116350534
$row9[0,20] = 3.05427116350534
$row9[0,21] = 3.05427116350534
$row9[0,22] = 3.05427116350534
$row9[0,23] = 3.05427116350534
$row9[0,24] = 3.05427116350534
$row9[0,25] = 3.05427116350534
$row9[0,26] = 3.05427116350534
$row9[0,27] = 3.05427116350534
$row9[0,28] = 3.05427116350534
$row9[0,29] = 3.05427116350534
$row9[0,30] = 3.05427116350534
$row9[0,31] = 3.05427116350534
$row9[0,32] = 3.05427116350534
$row9[0,33] = 3.05427116350534
$row9[0,34] = 3.05427116350534
$row9[0,35] = 3.05427116350534
$row9[0,36] = 3.05427116350534
$row9[0,37] = 3.05427116350534
$row9[0,38] = 3.05427116350534
$row9[0,39] = 3.05427116350534
$row9[0,40] = 3.05427116350534
$ws.Range("L9:AZ9").Value = $row9

$row10 = New-Object 'object[,]' 1,39
$row10[0,0] = 2.450141597588917
$row10[0,1] = 2.858912235977829
$row10[0,2] = 3.38738696315446
$row10[0,3] = 2.9587404276884
$row10[0,4] = 2.441258738366514
$row10[0,5] = 2.305809238174006
$row10[0,6] = 2.305809238174006
$row10[0,7] = 2.305809238174006
$row10[0,8] = 2.305809238174006
$row10[0,9] = 2.305809238174006
$row10[0,10] = 2.305809238174006
$row10[0,11] = 2.305809238174006
$row10[0,12] = 2.305809238174006
$row10[0,13] = 2.305809238174006
$row10[0,14] = 2.305809238174006
$row10[0,15] = 2.305809238174006
$row10[0,16] = 2.305809238174006
$row10[0,17] = 2.305809238174006
$row10[0,18] = 2.305809238174006
$row10[0,19] = 2.305809238174006
$row10[0,20] = 2.305809238174006
$row10[0,21] = 2.305809238174006
$row10[0,22] = 2.305809238174006
$row10[0,23] = 2.305809238174006
$row10[0,24] = 2.305809238174006
$row10[0,25] = 2.305809238174006
$row10[0,26] = 2.305809238174006
$row10[0,27] = 2.305809238174006
$row10[0,28] = 2.305809238174006
$row10[0,29] = 2.305809238174006
$row10[0,30] = 2.305809238174006
$row10[0,31] = 2.305809238174006
$row10[0,32] = 2.305809238174006
$row10[0,33] = 2.305809238174006
$row10[0,34] = 2.305809238174006
$row10[0,35] = 2.305809238174006
$row10[0,36] = 2.305809238174006
$row10[0,37] = 2.305809238174006
$row10[0,38] = 2.305809238174006
$ws.Range("N10:AZ10").Value = $row10

$row11 = New-Object 'object[,]' 1,37
$row11[0,0] = 3.491848178733536
$row11[0,1] = 3.466954089033747
$row11[0,2] = 3.091878630346012
$row11[0,3] = 2.661643377950096
$row11[0,4] = 2.480871685520603
$row11[0,5] = 2.509111342826809
$row11[0,6] = 2.509111342826809
$row11[0,7] = 2.509111342826809
$row11[0,8] = 2.509111342826809
$row11[0,9] = 2.509111342826809
$row11[0,10] = 2.509111342826809
$row11[0,11] = 2.509111342826809
$row11[0,12] = 2.509111342826809
$row11[0,13] = 2.509111342826809
$row11[0,14] = 2.509111342826809
$row11[0,15] = 2.509111342826809
$row11[0,16] = 2.509111342826809
$row11[0,17] = 2.509111342826809
$row11[0,18] = 2.509111342826809
$row11[0,19] = 2.509111342826809
$row11[0,20] = 2.509111342826809
$row11[0,21] = 2.509111342826809
$row11[0,22] = 2.509111342826809
$row11[0,23] = 2.509111342826809
$row11[0,24] = 2.509111342826809
$row11[0,25] = 2.509111342826809
$row11[0,26] = 2.509111342826809
$row11[0,27] = 2.509111342826809
$row11[0,28] = 2.509111342826809
$row11[0,29] = 2.509111342826809
$row11[0,30] = 2.509111342826809
$row11[0,31] = 2.509111342826809
$row11[0,32] = 2.509111342826809
$row11[0,33] = 2.509111342826809
$row11[0,34] = 2.509111342826809
$row11[0,35] = 2.509111342826809
$row11[0,36] = 2.509111342826809
$ws.Range("P11:AZ11").Value = $row11

$row12 = New-Object 'object[,]' 1,35
$row12[0,0] = 3.136162599657255
$row12[0,1] = 2.806147177668961
$row12[0,2] = 2.296583397191387
$row12[0,3] = 2.433248629349549
$row12[0,4] = 2.83347664679956
$row12[0,5] = 3.071095202329288
$row12[0,6] = 3.248721852957415
$row12[0,7] = 3.296731496509198
$row12[0,8] = 3.296731496509198
$row12[0,9] = 3.296731496509198
$row12[0,10] = 3.296731496509198
$row12[0,11] = 3.296731496509198
$row12[0,12] = 3.296731496509198
$row12[0,13] = 3.296731496509198
$row12[0,14] = 3.296731496509198
$row12[0,15] = 3.296731496509198
$row12[0,16] = 3.296731496509198
$row12[0,17] = 3.296731496509198
$row12[0,18] = 3.296731496509198
$row12[0,19] = 3.296731496509198
$row12[0,20] = 3.296731496509198
$row12[0,21] = 3.296731496509198
$row12[0,22] = 3.296731496509198
$row12[0,23] = 3.296731496509198
$row12[0,24] = 3.296731496509198
$row12[0,25] = 3.296731496509198
$row12[0,26] = 3.296731496509198
$row12[0,27] = 3.296731496509198
$row12[0,28] = 3.296731496509198
$row12[0,29] = 3.296731496509198
$row12[0,30] = 3.296731496509198
$row12[0,31] = 3.296731496509198
$row12[0,32] = 3.296731496509198
$row12[0,33] = 3.296731496509198
$row12[0,34] = 3.296731496509198
$ws.Range("R12:AZ12").Value = $row12

$row13 = New-Object 'object[,]' 1,33
$row13[0,0] = 2.2749004473406
$row13[0,1] = 2.419254005578297
$row13[0,2] = 2.560065157976177
$row13[0,3] = 2.781797072072023
$row13[0,4] = 3.078223990352669
$row13[0,5] = 3.296423324101938
$row13[0,6] = 3.149270133134596
$row13[0,7] = 2.994116795316071
$row13[0,8] = 2.867378798220366
$row13[0,9] = 2.861315725866587
$row13[0,10] = 2.861315725866587
$row13[0,11] = 2.861315725866587
$row13[0,12] = 2.861315725866587
$row13[0,13] = 2.861315725866587
$row13[0,14] = 2.861315725866587
$row13[0,15] = 2.861315725866587
$row13[0,16] = 2.861315725866587
$row13[0,17] = 2.861315725866587
$row13[0,18] = 2.861315725866587
$row13[0,19] = 2.861315725866587
$row13[0,20] = 2.861315725866587
$row13[0,21] = 2.861315725866587
$row13[0,22] = 2.861315725866587
$row13[0,23] = 2.861315725866587
$row13[0,24] = 2.861315725866587
$row13[0,25] = 2.861315725866587
$row13[0,26] = 2.861315725866587
$row13[0,27] = 2.861315725866587
$row13[0,28] = 2.861315725866587
$row13[0,29] = 2.861315725866587
$row13[0,30] = 2.861315725866587
$row13[0,31] = 2.861315725866587
$row13[0,32] = 2.861315725866587
$ws.Range("T13:AZ13").Value = $row13

$row14 = New-Object 'object[,]' 1,30
$row14[0,0] = 2.724246191199065
$row14[0,1] = 2.886418298927351
$row14[0,2] = 3.008574382540607
$row14[0,3] = 3.009352983329028
$row14[0,4] = 2.939737488252936
$row14[0,5] = 2.769017518462746
$row14[0,6] = 2.763966172716947
$row14[0,7] = 2.689501145820206
$row14[0,8] = 2.671604274379558
$row14[0,9] = 1.790319754067715
$row14[0,10] = 1.790319754067715
$row14[0,11] = 1.790319754067715
$row14[0,12] = 1.790319754067715
$row14[0,13] = 1.790319754067715
$row14[0,14] = 1.790319754067715
$row14[0,15] = 1.790319754067715
$row14[0,16] = 1.790319754067715
$row14[0,17] = 1.790319754067715
$row14[0,18] = 1.790319754067715
$row14[0,19] = 1.790319754067715
$row14[0,20] = 1.790319754067715
$row14[0,21] = 1.790319754067715
$row14[0,22] = 1.790319754067715
$row14[0,23] = 1.790319754067715
$row14[0,24] = 1.790319754067715
$row14[0,25] = 1.790319754067715
$row14[0,26] = 1.790319754067715
$row14[0,27] = 1.790319754067715
$row14[0,28] = 1.790319754067715
$row14[0,29] = 1.790319754067715
$ws.Range("W14:AZ14").Value = $row14

$row15 = New-Object 'object[,]' 1,26
$row15[0,0] = 2.988052171464251
$row15[0,1] = 2.927264293158816
$row15[0,2] = 2.947863484892133
$row15[0,3] = 2.967409274751098
$row15[0,4] = 2.997455747043043
$row15[0,5] = 1.627017245406992
$row15[0,6] = 1.856930494010856
$row15[0,7] = 1.980033360076905
$row15[0,8] = 2.08524086077817
$row15[0,9] = 2.339531676162721
$row15[0,10] = 2.339531676162721
$row15[0,11] = 2.339531676162721
$row15[0,12] = 2.339531676162721
$row15[0,13] = 2.339531676162721
$row15[0,14] = 2.339531676162721
$row15[0,15] = 2.339531676162721
$row15[0,16] = 2.339531676162721
$row15[0,17] = 2.339531676162721
$row15[0,18] = 2.339531676162721
$row15[0,19] = 2.339531676162721
$row15[0,20] = 2.339531676162721
$row15[0,21] = 2.339531676162721
$row15[0,22] = 2.339531676162721
$row15[0,23] = 2.339531676162721
$row15[0,24] = 2.339531676162721
$row15[0,25] = 2.339531676162721
$ws.Range("AA15:AZ15").Value = $row15

$row16 = New-Object 'object[,]' 1,22
$row16[0,0] = 3.008439268567842
$row16[0,1] = 2.174798403591915
$row16[0,2] = 2.092911340281423
$row16[0,3] = 2.063021041451907
$row16[0,4] = 2.197771900625956
$row16[0,5] = 3.195599391913406
$row16[0,6] = 4.270817433327112
$row16[0,7] = 4.865769161659883
$row16[0,8] = 4.939003803830477
$row16[0,9] = 4.834496776263886
$row16[0,10] = 4.834496776263886
$row16[0,11] = 4.834496776263886
$row16[0,12] = 4.834496776263886
$row16[0,13] = 4.834496776263886
$row16[0,14] = 4.834496776263886
$row16[0,15] = 4.834496776263886
$row16[0,16] = 4.834496776263886
$row16[0,17] = 4.834496776263886
$row16[0,18] = 4.834496776263886
$row16[0,19] = 4.834496776263886
$row16[0,20] = 4.834496776263886
$row16[0,21] = 4.834496776263886
$ws.Range("AE16:AZ16").Value = $row16

$row17 = New-Object 'object[,]' 1,19
$row17[0,0] = 2.012391101645061
$row17[0,1] = 2.053213017515065
$row17[0,2] = 2.441206385516637
$row17[0,3] = 2.849406056739201
$row17[0,4] = 3.291462037299842
$row17[0,5] = 3.440178795466697
$row17[0,6] = 3.06038938938058
$row17[0,7] = 3.604316462518464
$row17[0,8] = 3.332544669973525
$row17[0,9] = 2.93530792557688
$row17[0,10] = 2.798216547494237
$row17[0,11] = 2.798216547494237
$row17[0,12] = 2.798216547494237
$row17[0,13] = 2.798216547494237
$row17[0,14] = 2.798216547494237
$row17[0,15] = 2.798216547494237
$row17[0,16] = 2.798216547494237
$row17[0,17] = 2.798216547494237
$row17[0,18] = 2.798216547494237
$ws.Range("AH17:AZ17").Value = $row17

$row18 = New-Object 'object[,]' 1,15
$row18[0,0] = 3.13459343156206
$row18[0,1] = 3.2380444610977
$row18[0,2] = 3.116636734573786
$row18[0,3] = 3.937364994846959
$row18[0,4] = 3.789179157493971
$row18[0,5] = 3.215749572764803
$row18[0,6] = 2.721520966738655
$row18[0,7] = 2.098908173995873
$row18[0,8] = 1.888626610265987
$row18[0,9] = 1.635353376270698
$row18[0,10] = 1.530879676868468
$row18[0,11] = 1.530879676868468
$row18[0,12] = 1.530879676868468
$row18[0,13] = 1.530879676868468
$row18[0,14] = 1.530879676868468
$ws.Range("AL18:AZ18").Value = $row18

$row19 = New-Object 'object[,]' 1,11
$row19[0,0] = 3.87020515078067
$row19[0,1] = 3.630162063286146
$row19[0,2] = 3.491475308018321
$row19[0,3] = 3.330923984031142
$row19[0,4] = 2.949781091571957
$row19[0,5] = 2.328770194687713
$row19[0,6] = 1.713178787950698
$row19[0,7] = 1.874466487556892
$row19[0,8] = 1.966591496003445
$row19[0,9] = 1.984020855913604
$row19[0,10] = 2.060859685319461
$ws.Range("AP19:AZ19").Value = $row19

$row20 = New-Object 'object[,]' 1,7
$row20[0,0] = 2.975546095003945
$row20[0,1] = 2.484849225038532
$row20[0,2] = 1.999725833525323
$row20[0,3] = 1.913049717010873
$row20[0,4] = 2.005435469818684
$row20[0,5] = 1.88544721086894
$row20[0,6] = 2.141985433296578
$ws.Range("AT20:AZ20").Value = $row20

$row21 = New-Object 'object[,]' 1,3
$row21[0,0] = 1.983812695141185
$row21[0,1] = 1.801217086776363
$row21[0,2] = 1.846918513329565
$ws.Range("AX21:AZ21").Value = $row21
